# Dicionario de dados atualizado
# Add new "Transportadora_MeioTransporte" junction-table section to the
# "Fornecedor" sheet and remove the duplicated "id_meio_transporte" row
# from the "Cotacao" sheet (it now lives in the new section instead).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Fornecedor sheet: append the Transportadora_MeioTransporte table
# ---------------------------------------------------------------------
$wsForn = $wb.Worksheets.Item("Fornecedor")

# Row 70: section title (copy formatting from an existing section title,
# e.g. the "MeioTransporte" title on row 65)
$wsForn.Range("B65:G65").Copy()
$wsForn.Range("B70:G70").PasteSpecial(-4122)
$wsForn.Range("B70").Value = "Transportadora_MeioTransporte"
$wsForn.Range("B70:G70").Merge()

# Row 71: column header row (copy formatting from row 66)
$wsForn.Range("B66:G66").Copy()
$wsForn.Range("B71:G71").PasteSpecial(-4122)

# Row 72: id_transportadora (PK-FK) - copy formatting from a similar data row
$wsForn.Range("B67:G67").Copy()
$wsForn.Range("B72:G72").PasteSpecial(-4122)
$wsForn.Range("B72").Value = "id_transportadora"
$wsForn.Range("C72").Value = 9
$wsForn.Range("D72").Value = "Inteiro"
$wsForn.Range("E72").Value = 999999999
$wsForn.Range("F72").Value = "PK - FK"
$wsForn.Range("G72").Value = "Id da transportadora"

# Row 73: id_meio_transporte (PK-FK)
$wsForn.Range("B67:G67").Copy()
$wsForn.Range("B73:G73").PasteSpecial(-4122)
$wsForn.Range("B73").Value = "id_meio_transporte"
$wsForn.Range("C73").Value = 9
$wsForn.Range("D73").Value = "Inteiro"
$wsForn.Range("E73").Value = 999999999
$wsForn.Range("F73").Value = "PK - FK"
$wsForn.Range("G73").Value = "Id do meio de transporte"

$wsForn.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Cotacao sheet: remove the now-duplicated id_meio_transporte row
# ---------------------------------------------------------------------
$wsCot = $wb.Worksheets.Item("Cotacao")
$wsCot.Rows("24").Delete()

# ---------------------------------------------------------------------
# 3) Leave "Fornecedor" as the active sheet (matches the saved file)
# ---------------------------------------------------------------------
$wsForn.Activate()
